$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("B3").Value = "['MEC-3B-Trat. Termicos', -, -, -]"
$ws.Range("B4").Value = "['MEC-3B-Trat. Termicos', -, -, -]"
$ws.Range("B6").Value = "['MEC-3B-Trat. Termicos', -, -, -]"
$ws.Range("E6").Value = "-"
$ws.Range("B7").Value = "['MEC-3B-Trat. Termicos', -, -, -]"
$ws.Range("E7").Value = "-"
$ws.Range("B8").Value = "-"
$ws.Range("E10").Value = "-"
$ws.Range("E11").Value = "[-, -, 'MEC-3A-Tec. Soldagem', -]"
$ws.Range("E12").Value = "[-, -, 'MEC-3A-Tec. Soldagem', -]"
$ws.Range("E14").Value = "[-, -, 'MEC-3A-Tec. Soldagem', -]"
$ws.Range("F14").Value = "-"
$ws.Range("E15").Value = "[-, -, 'MEC-3A-Tec. Soldagem', -]"
$ws.Range("F15").Value = "-"
$ws.Range("E16").Value = "-"
